$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "nie"
$ws.Range("B2").Value = 5
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "696273167"
$ws.Range("D2").Value = "Samugn"
$ws.Range("E2").Value = "13/11/2023"
$ws.Range("F2").Value = "opsi`t"
$ws.Range("G2").Value = "nas"

# Row 3
$ws.Range("B3").Value = 6
$ws.Range("C3").Value = "test"
$ws.Range("D3").Value = "ddd"
$ws.Range("E3").Value = "13/11/2023"
$ws.Range("F3").Value = "oooo`t"
$ws.Range("G3").Value = "nnnnn"

# Row 4
$ws.Range("B4").Value = 7
$ws.Range("C4").Value = "678 098"
$ws.Range("D4").Value = "To ejst mode"
$ws.Range("E4").Value = "13/11/2023"
$ws.Range("F4").Value = "dasddasdjkksdfjsdahfaskfslkfshdjflaflsf`tsdkjf"
$ws.Range("G4").Value = "ndfmnfnlkjlkjldsjklskjldsfljkdsfljkdfljk"
